$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the original row 2 and row 3 values before overwriting, so we can swap them.
# Use Value2 (not Value) to avoid COM Variant wrapper objects being stringified.
$origA2 = $ws.Range("A2").Value2
$origD2 = $ws.Range("D2").Value2
$origE2 = $ws.Range("E2").Value2
$origF2 = $ws.Range("F2").Value2
$origG2 = $ws.Range("G2").Value2
$origH2 = $ws.Range("H2").Value2
$origQ2 = $ws.Range("Q2").Value2
$origR2 = $ws.Range("R2").Value2
$origZ2 = $ws.Range("Z2").Value2
$origAB2 = $ws.Range("AB2").Value2

$origA3 = $ws.Range("A3").Value2
$origD3 = $ws.Range("D3").Value2
$origE3 = $ws.Range("E3").Value2
$origF3 = $ws.Range("F3").Value2
$origG3 = $ws.Range("G3").Value2
$origH3 = $ws.Range("H3").Value2
$origQ3 = $ws.Range("Q3").Value2
$origR3 = $ws.Range("R3").Value2
$origZ3 = $ws.Range("Z3").Value2
$origAB3 = $ws.Range("AB3").Value2

# Row 2 becomes what was row 3 (content swap), except B2 gets a fresh value.
$ws.Range("A2").Value2 = $origA3
$ws.Range("B2").Value2 = 89553
$ws.Range("D2").Value2 = $origD3
$ws.Range("E2").Value2 = $origE3
$ws.Range("F2").Value2 = $origF3
$ws.Range("G2").Value2 = $origG3
$ws.Range("H2").Value2 = $origH3
$ws.Range("Q2").Value2 = $origQ3
$ws.Range("R2").Value2 = $origR3
$ws.Range("Z2").Value2 = "12:44"
$ws.Range("AB2").Value2 = "12:44"

# Row 3 becomes what was row 2 (content swap), except B3 gets a fresh value.
$ws.Range("A3").Value2 = $origA2
$ws.Range("B3").Value2 = 95921
$ws.Range("D3").Value2 = $origD2
$ws.Range("E3").Value2 = $origE2
$ws.Range("F3").Value2 = $origF2
$ws.Range("G3").Value2 = $origG2
$ws.Range("H3").Value2 = $origH2
$ws.Range("Q3").Value2 = $origQ2
$ws.Range("R3").Value2 = $origR2
$ws.Range("Z3").Value2 = "12:41"
$ws.Range("AB3").Value2 = "12:41"
